# Remove the "Note:" paragraph (explaining the assumed identity of a data
# scientist) along with the blank paragraph that immediately follows it.
# The blank paragraph before "Note:" is left untouched.

$d = $word.ActiveDocument

# Locate the paragraph containing the "Note:" text using Find (wildcards on,
# so we don't depend on the exact ellipsis text matching character-for-character).
$rng = $d.Content
$found = $rng.Find.Execute(
    "Note:*identity of a data scientist*ankle braces, etc).",
    $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Note:' paragraph to remove."
}

# Figure out which paragraph (by index) the match starts in.
$startIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $rng.Start) {
        $startIdx = $i
        break
    }
}

if ($startIdx -eq -1) {
    throw "Could not resolve the paragraph index for the found range."
}

# The "Note:" paragraph plus the following (blank) paragraph get removed.
$noteParagraph  = $d.Paragraphs.Item($startIdx)
$blankParagraph = $d.Paragraphs.Item($startIdx + 1)
$afterParagraph = $d.Paragraphs.Item($startIdx + 2)

$deleteRange = $d.Range($noteParagraph.Range.Start, $afterParagraph.Range.Start)
$deleteRange.Delete()
